# Generate Report for Handback
# Updates the Overview status text, records handback xliff target/file/datetime
# info for the zh-cn and de-de sheets, adds "Latest Target File" hyperlinks,
# and widens the columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$aUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f826c960ad716232de918ed72c3568cf89b121b5/e2e/a.md"
$bUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f826c960ad716232de918ed72c3568cf89b121b5/e2e/b.md"

# --- Overview sheet: status is now "Handed back: in sync with en-US" ---
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.16
$overview.Columns.Item(6).ColumnWidth = 29.16

# --- zh-cn sheet: record the generated handback xliff + target file + datetime ---
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-04 16:41:55"
$zhcn.Range("K3").Value = "2016-09-04 16:41:55"

$zhcn.Columns.Item(3).ColumnWidth = 29.16
$zhcn.Columns.Item(10).ColumnWidth = 39.17

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $aUrl, [Type]::Missing, [Type]::Missing, "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $aUrl, [Type]::Missing, [Type]::Missing, "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $bUrl, [Type]::Missing, [Type]::Missing, "b.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $aUrl, [Type]::Missing, [Type]::Missing, "a.md")

# --- de-de sheet: record the generated handback xliff + target file + datetime ---
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-09-04 16:42:09"
$dede.Range("K3").Value = "2016-09-04 16:42:09"

$dede.Columns.Item(3).ColumnWidth = 29.16
$dede.Columns.Item(10).ColumnWidth = 39.17

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $aUrl, [Type]::Missing, [Type]::Missing, "a.md")
$dede.Hyperlinks.Add($dede.Range("I2"), $aUrl, [Type]::Missing, [Type]::Missing, "a.md")
$dede.Hyperlinks.Add($dede.Range("A3"), $bUrl, [Type]::Missing, [Type]::Missing, "b.md")
$dede.Hyperlinks.Add($dede.Range("I3"), $aUrl, [Type]::Missing, [Type]::Missing, "a.md")
